# Refresh cryptocurrency price/volume snapshot (GitHub Actions update, 2024-01-05).
# Rebuilds the table rows that moved between scrapes, including two rows whose
# ranking order swapped (Kaspa/ARBITRUM and MultiversX/Stacks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$text) {
    # Writing via COM auto-detects numeric-looking strings as numbers, which would
    # drop significant trailing zeros (e.g. "14.10" -> 14.1). Prefix those with an
    # apostrophe (the normal Excel 'treat as text' marker) to keep the literal text.
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $ws.Range($cellRef).Value = "'" + $text
    } else {
        $ws.Range($cellRef).Value = $text
    }
}

# Row 2
Set-TextValue "D2" '43.993.90'
Set-TextValue "E2" '  +1.39%  '

# Row 3
Set-TextValue "D3" '2.241.63'
Set-TextValue "E3" '  +0.13%  '

# Row 4
Set-TextValue "E4" '  +0.18%  '

# Row 5
Set-TextValue "D5" '315.94'
Set-TextValue "E5" '  -1.00%  '

# Row 6
Set-TextValue "D6" '100.42'
Set-TextValue "E6" '  +0.24%  '

# Row 7
Set-TextValue "D7" '0.572'
Set-TextValue "E7" '  -2.04%  '

# Row 8
Set-TextValue "E8" '  +0.18%  '

# Row 9
Set-TextValue "D9" '0.540'
Set-TextValue "E9" '  -4.14%  '

# Row 10
Set-TextValue "D10" '36.57'
Set-TextValue "E10" '  -1.99%  '

# Row 11
Set-TextValue "D11" '0.0824'
Set-TextValue "E11" '  -1.10%  '

# Row 12
Set-TextValue "D12" '7.43'
Set-TextValue "E12" '  -3.84%  '

# Row 13
Set-TextValue "E13" '  -2.53%  '

# Row 14
Set-TextValue "D14" '2.584.04'
Set-TextValue "E14" '  +0.20%  '

# Row 15
Set-TextValue "D15" '0.845'
Set-TextValue "E15" '  -2.58%  '

# Row 16
Set-TextValue "D16" '2.247.94'
Set-TextValue "E16" '  +0.56%  '

# Row 17
Set-TextValue "D17" '14.10'
Set-TextValue "E17" '  -1.58%  '

# Row 18
Set-TextValue "D18" '43.923.19'
Set-TextValue "E18" '  +1.29%  '

# Row 19
Set-TextValue "D19" '13.03'
Set-TextValue "E19" '  -8.80%  '

# Row 20
Set-TextValue "D20" '0.0₃0967'
Set-TextValue "E20" '  -1.11%  '

# Row 21
Set-TextValue "D21" '6.40'
Set-TextValue "E21" '  -3.57%  '

# Row 22
Set-TextValue "D22" '65.16'
Set-TextValue "E22" '  -0.64%  '

# Row 23
Set-TextValue "D23" '3.06'
Set-TextValue "E23" '  -4.24%  '

# Row 24
Set-TextValue "D24" '234.64'
Set-TextValue "E24" '  -0.97%  '

# Row 25
Set-TextValue "E25" '  -5.85%  '

# Row 26
Set-TextValue "E26" '  +0.16%  '

# Row 27
Set-TextValue "D27" '10.34'
Set-TextValue "E27" '  +2.57%  '

# Row 28
Set-TextValue "E28" '  -0.58%  '

# Row 29
Set-TextValue "D29" '36.85'
Set-TextValue "E29" '  +0.76%  '

# Row 30
Set-TextValue "D30" '6.08'
Set-TextValue "E30" '  -4.97%  '

# Row 31
Set-TextValue "D31" '159.00'
Set-TextValue "E31" '  -0.07%  '

# Row 32
Set-TextValue "D32" '20.04'
Set-TextValue "E32" '  -1.44%  '

# Row 33
Set-TextValue "D33" '0.0840'
Set-TextValue "E33" '  -3.95%  '

# Row 34
Set-TextValue "D34" '2.69'
Set-TextValue "E34" '  -1.01%  '

# Row 35
Set-TextValue "D35" '3.16'
Set-TextValue "E35" '  -1.93%  '

# Row 36
Set-TextValue "B36" 'Kaspa'
Set-TextValue "C36" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D36" '0.111'
Set-TextValue "E36" '  +6.24%  '

# Row 37
Set-TextValue "B37" 'ARBITRUM'
Set-TextValue "C37" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D37" '1.92'
Set-TextValue "E37" '  +1.12%  '

# Row 38
Set-TextValue "E38" '  -2.47%  '

# Row 39
Set-TextValue "D39" '15.97'
Set-TextValue "E39" '  +9.78%  '

# Row 40
Set-TextValue "D40" '3.64'
Set-TextValue "E40" '  -2.19%  '

# Row 41
Set-TextValue "D41" '4.08'
Set-TextValue "E41" '  -6.80%  '

# Row 42
Set-TextValue "D42" '0.0312'
Set-TextValue "E42" '  -3.46%  '

# Row 43
Set-TextValue "E43" '  +0.14%  '

# Row 44
Set-TextValue "D44" '1.734.19'
Set-TextValue "E44" '  -4.90%  '

# Row 45
Set-TextValue "D45" '0.196'
Set-TextValue "E45" '  -4.36%  '

# Row 46
Set-TextValue "D46" '81.24'
Set-TextValue "E46" '  -3.59%  '

# Row 47
Set-TextValue "D47" '73.68'
Set-TextValue "E47" '  -1.60%  '

# Row 48
Set-TextValue "D48" '5.13'
Set-TextValue "E48" '  -3.26%  '

# Row 49
Set-TextValue "D49" '102.01'
Set-TextValue "E49" '  -1.15%  '

# Row 50
Set-TextValue "B50" 'MultiversX'
Set-TextValue "C50" 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue "D50" '57.24'
Set-TextValue "E50" '  -2.43%  '

# Row 51
Set-TextValue "B51" 'Stacks'
Set-TextValue "C51" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D51" '1.64'
Set-TextValue "E51" '  +0.30%  '

